# Lab30421 - update lab test grades (R/S columns), recompute the
# dependent "final mark" (T) and "situation" (U) columns as plain
# values (the old formulas are replaced by their computed results),
# and refresh the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R3").Value = 10
$ws.Range("S3").Value = 4.3
$ws.Range("T3").Value = 7.35
$ws.Range("U3").Value = "fail"

$ws.Range("R4").Value = 6.5
$ws.Range("S4").Value = 3
$ws.Range("T4").Value = 4.9000000000000004
$ws.Range("U4").Value = "fail"

$ws.Range("R5").Value = 10
$ws.Range("S5").Value = 8.1999999999999993
$ws.Range("T5").Value = 9.6
$ws.Range("U5").Value = "pass"

$ws.Range("R6").Value = 7
$ws.Range("S6").Value = 3.2
$ws.Range("T6").Value = 5.4
$ws.Range("U6").Value = "fail"

$ws.Range("R7").Value = 9.5
$ws.Range("S7").Value = 7.6
$ws.Range("T7").Value = 8.8000000000000007
$ws.Range("U7").Value = "pass"

$ws.Range("R8").Value = 4
$ws.Range("S8").Value = 2.8
$ws.Range("T8").Value = 3.5
$ws.Range("U8").Value = "fail"

$ws.Range("R9").Value = 8.75
$ws.Range("S9").Value = 8.9
$ws.Range("T9").Value = 8.9749999999999996
$ws.Range("U9").Value = "pass"

$ws.Range("R10").Value = 6
$ws.Range("S10").Value = 3.1
$ws.Range("T10").Value = 4.68
$ws.Range("U10").Value = "fail"

$ws.Range("R11").Value = 6.5
$ws.Range("S11").Value = 3.1
$ws.Range("T11").Value = 4.88
$ws.Range("U11").Value = "fail"

$ws.Range("R12").Value = 8.5
$ws.Range("S12").Value = 4.0999999999999996
$ws.Range("T12").Value = 6.38
$ws.Range("U12").Value = "fail"

$ws.Range("R13").Value = 10
$ws.Range("S13").Value = 7.1
$ws.Range("T13").Value = 8.73
$ws.Range("U13").Value = "pass"

$ws.Range("R14").Value = 9.25
$ws.Range("S14").Value = 9.1999999999999993
$ws.Range("T14").Value = 9.4049999999999994
$ws.Range("U14").Value = "pass"

$ws.Range("R15").Value = 9.5
$ws.Range("S15").Value = 7.9
$ws.Range("T15").Value = 8.89
$ws.Range("U15").Value = "pass"

$ws.Range("R16").Value = 10
$ws.Range("S16").Value = 5.8
$ws.Range("T16").Value = 8
$ws.Range("U16").Value = "pass"

$ws.Range("R17").Value = 9
$ws.Range("S17").Value = 7.4
$ws.Range("T17").Value = 8.6999999999999993
$ws.Range("U17").Value = "pass"

$ws.Range("R18").Value = 1
$ws.Range("S18").Value = 1
$ws.Range("T18").Value = 1
$ws.Range("U18").Value = "fail"

$ws.Range("R19").Value = 10
$ws.Range("S19").Value = 8.3000000000000007
$ws.Range("T19").Value = 9.6
$ws.Range("U19").Value = "pass"

$ws.Range("R20").Value = 9.5
$ws.Range("S20").Value = 8.1999999999999993
$ws.Range("T20").Value = 8.85
$ws.Range("U20").Value = "pass"

$ws.Range("R21").Value = 10
$ws.Range("S21").Value = 8.4
$ws.Range("T21").Value = 10
$ws.Range("U21").Value = "pass"

$ws.Range("R22").Value = 5.5
$ws.Range("S22").Value = 7
$ws.Range("T22").Value = 6.37
$ws.Range("U22").Value = "pass"

$ws.Range("R23").Value = 1
$ws.Range("S23").Value = 1
$ws.Range("T23").Value = 1
$ws.Range("U23").Value = "fail"

$ws.Range("R24").Value = 6
$ws.Range("S24").Value = 5.0999999999999996
$ws.Range("T24").Value = 5.7
$ws.Range("U24").Value = "pass"

$ws.Range("R25").Value = 10
$ws.Range("S25").Value = 9.6
$ws.Range("T25").Value = 10
$ws.Range("U25").Value = "pass"

# Match the author's final on-screen selection
$ws.Range("Q2:U25").Select()
